# "market maker seems okay on 1017" - update market-maker (J column) trade
# inputs on the Diversification sheet and add the two supporting cells that
# compare the hedge payoff (L20) against the new benchmark (N18).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diversification")

# Market-maker proposed-trade inputs (column J, rows 8-12)
$ws.Range("J8").Value  = 0.04
$ws.Range("J9").Value  = -10
$ws.Range("J10").Value = 10
$ws.Range("J11").Value = 10
$ws.Range("J12").Value = 9

# New benchmark figure and its comparison against the hedge payoff (L20)
$ws.Range("N18").Value = 27.8525
$ws.Range("O19").Formula = "=L20-N18"

# The source rows auto-fit to a shorter height once recalculated
$ws.Rows.Item(8).RowHeight = 19.7
$ws.Rows.Item(9).RowHeight = 19.7
$ws.Rows.Item(10).RowHeight = 19.7
$ws.Rows.Item(12).RowHeight = 19.7

# Leave the cursor where the author left it
$ws.Range("J11").Select()
